# Handback status report refresh: a later handback run completed for
# 9c45f0c6-5b22-4ba6-9214-cdd1bb6eaf0b, so its timestamps/priority move
# forward. Note: rows for 9c45f0c6-... and bb3228b0-... happen to share the
# exact same datetime text in several columns (Excel stores identical
# strings once), so updating the 9c45f0c6 row's value also changes the
# bb3228b0 row's displayed value in those columns - that is expected and is
# reproduced here explicitly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" column (G) for rows 3 (9c45f0c6-...) and 4
# (bb3228b0-...) - both currently read 2016-08-19 14:15:19.
$wsOverview.Range("G3").Value = "2016-08-19 14:16:12"
$wsOverview.Range("G4").Value = "2016-08-19 14:16:12"

# --- zh-cn sheet ---
# "Priority" column (E) for rows 3 and 4 - both currently read "ht".
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# "Correspond Handoff Datetime" column (H) for rows 3 and 4 - both currently
# read 2016-08-19 14:15:03.
$wsZhCn.Range("H3").Value = "2016-08-19 14:15:58"
$wsZhCn.Range("H4").Value = "2016-08-19 14:15:58"

# "Correspond Handback DateTime" column (K) for rows 3 and 4 - both
# currently read 2016-08-19 14:15:32.
$wsZhCn.Range("K3").Value = "2016-08-19 14:16:30"
$wsZhCn.Range("K4").Value = "2016-08-19 14:16:30"

# --- de-de sheet ---
# "Correspond Handoff Datetime" column (H) for rows 3 and 4 - both currently
# read 2016-08-19 14:15:19.
$wsDeDe.Range("H3").Value = "2016-08-19 14:16:12"
$wsDeDe.Range("H4").Value = "2016-08-19 14:16:12"

# "Priority" column (E) for rows 3 and 4 - both currently read "ht".
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# "Correspond Handback DateTime" column (K) for rows 3 and 4 - both
# currently read 2016-08-19 14:15:45.
$wsDeDe.Range("K3").Value = "2016-08-19 14:16:37"
$wsDeDe.Range("K4").Value = "2016-08-19 14:16:37"
